# "Drop in RMI script results for 3.0"
#
# This workbook previously derived Texas-specific elasticities on a
# "Texas Notes" sheet (hand-measured from pixel areas in a source graphic)
# and averaged them with the national "Calculations" numbers on the
# "EoDSDwSP" sheet. The update drops the Texas-specific workaround and
# points EoDSDwSP straight at the national "Calculations" results.

$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

$wsEoDSDwSP   = $wb.Worksheets.Item("EoDSDwSP")
$wsTexasNotes = $wb.Worksheets.Item("Texas Notes")

# 1. Re-point EoDSDwSP's elasticity formulas at Calculations (national
#    numbers) instead of the Texas Notes averages, before the Texas Notes
#    sheet is removed.
$wsEoDSDwSP.Range("B2").Formula = "=Calculations!B9"
$wsEoDSDwSP.Range("B4").Formula = "=Calculations!B10"

# 2. Drop the now-unused "Texas Notes" sheet and its workings entirely.
$wsTexasNotes.Delete()

# Re-fetch worksheet references: deleting a sheet invalidates the
# references the remaining worksheets were bound to beforehand.
$wsAbout    = $wb.Worksheets.Item("About")
$wsCalc     = $wb.Worksheets.Item("Calculations")
$wsEoDSDwSP = $wb.Worksheets.Item("EoDSDwSP")

# 3. The "About" sheet's source link is no longer a live hyperlink -
#    remove it (keep the cell text) and restore the cell to normal style.
$wsAbout.Hyperlinks.Delete()
$wsAbout.Range("B6").Style = "Normal"

# 4. Update the cursor/selection state left in each sheet and leave the
#    "About" sheet as the active tab.
$wsCalc.Activate()
$wsCalc.Range("A1").Select()

$wsEoDSDwSP.Activate()
$wsEoDSDwSP.Range("B2").Select()

$wsAbout.Activate()
$wsAbout.Range("A12").Select()
